$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated symbol list values (text-preserving) from coinranking.com scrape
$updates = @(
    @{Cell='D2'; Value='328.78'},
    @{Cell='E2'; Value='0.03%'},
    @{Cell='D3'; Value='44.37'},
    @{Cell='E3'; Value='-0.07%'},
    @{Cell='D4'; Value='5.506'},
    @{Cell='E4'; Value='-1.59%'},
    @{Cell='D5'; Value='0.08065'},
    @{Cell='E5'; Value='-0.27%'},
    @{Cell='D6'; Value='2.055'},
    @{Cell='E6'; Value='1.07%'},
    @{Cell='D7'; Value='0.9614'},
    @{Cell='E7'; Value='1.09%'},
    @{Cell='D8'; Value='0.1122'},
    @{Cell='E8'; Value='-3.95%'},
    @{Cell='E9'; Value='1.51%'},
    @{Cell='D10'; Value='10.19'},
    @{Cell='E10'; Value='-0.58%'},
    @{Cell='D11'; Value='0.09911'},
    @{Cell='E11'; Value='2.15%'},
    @{Cell='D12'; Value='0.04719'},
    @{Cell='E12'; Value='2.70%'},
    @{Cell='D13'; Value='0.1064'},
    @{Cell='E13'; Value='-0.41%'},
    @{Cell='E14'; Value='-1.96%'},
    @{Cell='D15'; Value='0.04101'},
    @{Cell='E15'; Value='-2.43%'},
    @{Cell='D16'; Value='0.006132'},
    @{Cell='E16'; Value='2.06%'},
    @{Cell='B17'; Value='HotbitToken'},
    @{Cell='C17'; Value='https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'},
    @{Cell='D17'; Value='0.004335'},
    @{Cell='E17'; Value='0.49%'},
    @{Cell='B18'; Value='LEO'},
    @{Cell='C18'; Value='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'},
    @{Cell='D18'; Value='3.338'},
    @{Cell='E18'; Value='-0.88%'},
    @{Cell='B19'; Value='GateToken'},
    @{Cell='C19'; Value='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'},
    @{Cell='D19'; Value='4.430'},
    @{Cell='E19'; Value='2.98%'},
    @{Cell='B20'; Value='BTSEToken'},
    @{Cell='C20'; Value='https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'},
    @{Cell='D20'; Value='2.621'},
    @{Cell='E20'; Value='4.12%'},
    @{Cell='B21'; Value='BitpandaEcosystemToken'},
    @{Cell='C21'; Value='https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'},
    @{Cell='D21'; Value='0.3313'},
    @{Cell='E21'; Value='-4.86%'},
    @{Cell='B22'; Value='ProBitToken'},
    @{Cell='C22'; Value='https://coinranking.com/coin/lQP4d6T2+probittoken-prob'},
    @{Cell='D22'; Value='0.1394'},
    @{Cell='E22'; Value='-1.14%'},
    @{Cell='B23'; Value='ZBToken'},
    @{Cell='C23'; Value='https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'},
    @{Cell='D23'; Value='0.2581'},
    @{Cell='E23'; Value='2.97%'},
    @{Cell='B24'; Value='BitKan'},
    @{Cell='C24'; Value='https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'},
    @{Cell='D24'; Value='0.001313'},
    @{Cell='E24'; Value='5.29%'},
    @{Cell='D25'; Value='0.0001285'},
    @{Cell='E25'; Value='7.94%'},
    @{Cell='D26'; Value='0.0003753'},
    @{Cell='E26'; Value='-5.68%'},
    @{Cell='D38'; Value='0.02631'},
    @{Cell='E38'; Value='-1.15%'},
    @{Cell='D39'; Value='0.05613'},
    @{Cell='E39'; Value='1.26%'},
    @{Cell='D40'; Value='0.007637'},
    @{Cell='E40'; Value='0.72%'},
    @{Cell='E41'; Value='-0.35%'},
    @{Cell='D42'; Value='0.007419'},
    @{Cell='E42'; Value='-8.20%'},
    @{Cell='D43'; Value='0.001993'},
    @{Cell='E43'; Value='-1.19%'},
    @{Cell='D44'; Value='0.008732'},
    @{Cell='E44'; Value='3.99%'},
    @{Cell='D45'; Value='0.00007125'},
    @{Cell='E45'; Value='-0.62%'},
    @{Cell='D46'; Value='0.00000000753'},
    @{Cell='E46'; Value='0.21%'},
    @{Cell='D47'; Value='0.0005820'},
    @{Cell='E47'; Value='0.15%'},
    @{Cell='D48'; Value='0.002529'},
    @{Cell='E48'; Value='11.24%'},
    @{Cell='D49'; Value='0.003572'},
    @{Cell='E49'; Value='-20.67%'},
    @{Cell='E50'; Value='0.21%'},
    @{Cell='D51'; Value='0.0002007'},
    @{Cell='E51'; Value='0.21%'}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
